$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Preserve the existing text storage type (these cells hold numeric-looking
# strings, not real numbers) by temporarily switching to Text format while
# assigning, then resetting the style back to Normal so no stray number
# format is left attached to the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "289.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.456"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06378"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "6.577"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8282"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01426"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1687"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08811"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03656"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09188"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.730"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001659"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006192"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006299"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001073"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.782"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.406"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3360"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1263"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002711"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04835"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007144"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004508"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1115"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01180"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007040"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8015"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.006525"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00001904"
$ws.Range("D49").Style = "Normal"

# --- Column G (Hora) updates: all rows 2-51 go from "20" to "21" ---
$gRange = $ws.Range("G2:G51")
$gRange.NumberFormat = "@"
$ws.Range("G2").Value = "21"
$ws.Range("G3").Value = "21"
$ws.Range("G4").Value = "21"
$ws.Range("G5").Value = "21"
$ws.Range("G6").Value = "21"
$ws.Range("G7").Value = "21"
$ws.Range("G8").Value = "21"
$ws.Range("G9").Value = "21"
$ws.Range("G10").Value = "21"
$ws.Range("G11").Value = "21"
$ws.Range("G12").Value = "21"
$ws.Range("G13").Value = "21"
$ws.Range("G14").Value = "21"
$ws.Range("G15").Value = "21"
$ws.Range("G16").Value = "21"
$ws.Range("G17").Value = "21"
$ws.Range("G18").Value = "21"
$ws.Range("G19").Value = "21"
$ws.Range("G20").Value = "21"
$ws.Range("G21").Value = "21"
$ws.Range("G22").Value = "21"
$ws.Range("G23").Value = "21"
$ws.Range("G24").Value = "21"
$ws.Range("G25").Value = "21"
$ws.Range("G26").Value = "21"
$ws.Range("G27").Value = "21"
$ws.Range("G28").Value = "21"
$ws.Range("G29").Value = "21"
$ws.Range("G30").Value = "21"
$ws.Range("G31").Value = "21"
$ws.Range("G32").Value = "21"
$ws.Range("G33").Value = "21"
$ws.Range("G34").Value = "21"
$ws.Range("G35").Value = "21"
$ws.Range("G36").Value = "21"
$ws.Range("G37").Value = "21"
$ws.Range("G38").Value = "21"
$ws.Range("G39").Value = "21"
$ws.Range("G40").Value = "21"
$ws.Range("G41").Value = "21"
$ws.Range("G42").Value = "21"
$ws.Range("G43").Value = "21"
$ws.Range("G44").Value = "21"
$ws.Range("G45").Value = "21"
$ws.Range("G46").Value = "21"
$ws.Range("G47").Value = "21"
$ws.Range("G48").Value = "21"
$ws.Range("G49").Value = "21"
$ws.Range("G50").Value = "21"
$ws.Range("G51").Value = "21"
$gRange.Style = "Normal"
